$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# The rows for SingleUseId1 (row 4), SingleUseId4 (row 5) and SingleUseId8 (row 9)
# are no longer needed now that there are two analog dials (currentRpm /
# totalDistance) sharing a single config - delete those three rows, shifting
# the remaining rows up.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
